$d = $word.ActiveDocument

$replacements = @(
    @("19÷2=9, 1", "56÷5=11, 1"),
    @("43÷7=6, 1", "20÷8=2, 4"),
    @("13÷7=1, 6", "25÷8=3, 1"),
    @("66÷3=22, 0", "79÷4=19, 3"),
    @("34÷2=17, 0", "79÷5=15, 4"),
    @("87÷8=10, 7", "75÷5=15, 0"),
    @("55÷3=18, 1", "32÷7=4, 4"),
    @("10÷4=2, 2", "87÷4=21, 3"),
    @("55÷2=27, 1", "67÷9=7, 4"),
    @("66÷9=7, 3", "79÷5=15, 4"),
    @("10÷7=1, 3", "34÷8=4, 2"),
    @("49÷7=7, 0", "67÷3=22, 1"),
    @("35÷6=5, 5", "92÷3=30, 2"),
    @("95÷2=47, 1", "87÷3=29, 0"),
    @("70÷8=8, 6", "98÷9=10, 8"),
    @("23÷2=11, 1", "40÷8=5, 0"),
    @("71÷4=17, 3", "60÷5=12, 0"),
    @("33÷5=6, 3", "91÷7=13, 0"),
    @("96÷6=16, 0", "84÷5=16, 4"),
    @("25÷3=8, 1", "26÷5=5, 1"),
    @("85÷5=17, 0", "39÷8=4, 7"),
    @("24÷2=12, 0", "89÷5=17, 4"),
    @("49÷3=16, 1", "80÷9=8, 8"),
    @("42÷7=6, 0", "28÷6=4, 4"),
    @("12÷3=4, 0", "91÷5=18, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
